$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh Price (D) and Volume(1h) (E) columns with latest scraped values
$ws.Range("D2").Value = "28.435.96"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "1.863.03"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'324.63"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").Value = "'1.007"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "'0.4565"
$ws.Range("E7").Value = "  -1.68%  "
$ws.Range("D8").Value = "'0.3836"
$ws.Range("E8").Value = "  -0.93%  "
$ws.Range("D9").Value = "'0.07820"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").Value = "'0.9865"
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("D11").Value = "'21.53"
$ws.Range("E11").Value = "  -2.70%  "
$ws.Range("D12").Value = "1.883.17"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").Value = "'6.910"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").Value = "'5.633"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").Value = "'0.06972"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").Value = "'86.90"
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("D17").Value = "'1.008"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "'0.000009977"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").Value = "'16.65"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").Value = "'1.008"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "28.449.06"
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("D22").Value = "'5.257"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").Value = "'10.91"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("D24").Value = "'2.105"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "2.139.19"
$ws.Range("E25").Value = "  +3.32%  "
$ws.Range("D26").Value = "'153.05"
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("D27").Value = "'19.14"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").Value = "'5.679"
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("D29").Value = "'1.938"
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("D30").Value = "'117.83"
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("D31").Value = "'0.09276"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "'0.9092"
$ws.Range("E32").Value = "  -2.79%  "
$ws.Range("D33").Value = "'5.275"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").Value = "'1.319"
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("D35").Value = "'3.304"
$ws.Range("E35").Value = "  -0.95%  "
$ws.Range("D36").Value = "'0.05718"
$ws.Range("E36").Value = "  -2.01%  "
$ws.Range("D37").Value = "'1.138"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "'0.02058"
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("D39").Value = "'7.672"
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("D40").Value = "'0.5564"
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("D41").Value = "'0.1772"
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("D42").Value = "'9.651"
$ws.Range("E42").Value = "  -2.58%  "
$ws.Range("D43").Value = "'0.07098"
$ws.Range("E43").Value = "  -2.44%  "
$ws.Range("D44").Value = "'11.59"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D47").Value = "'1.112"
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("D48").Value = "'1.815"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("D49").Value = "'112.01"
$ws.Range("E49").Value = "  -1.78%  "
$ws.Range("D50").Value = "'2.412"
$ws.Range("E50").Value = "  +3.79%  "

# RenderToken moved above Decentraland in the ranking; update rows 45-46 accordingly
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'2.156"
$ws.Range("E45").Value = "  +0.98%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5236"
$ws.Range("E46").Value = "  -0.87%  "

# Volume(1h) refresh for the final row
$ws.Range("E51").Value = "  +0.08%  "
